# Generate Report for Handoff
#
# The localization status has moved from "Handed back: in sync with en-US"
# to "Ready for handoff", with a refreshed timestamp, the "ht" job-type code
# corrected to "mt", a refreshed target-generation timestamp on the zh-cn
# sheet, and a "stale handback" error message recorded for the 5c36ebef...
# file on both the zh-cn and de-de report sheets.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Overview")
$ws2 = $wb.Worksheets.Item("zh-cn")
$ws3 = $wb.Worksheets.Item("de-de")

# --- Overview sheet -------------------------------------------------------
# Status text + Latest HO Xliff Generate Date for both language rows.
$ws1.Range("E2:F3").Value = "Ready for handoff"
$ws1.Range("G2").Value = "2016-11-09 01:39:24"
$ws1.Range("G3").Value = "2016-11-09 01:39:24"

# Overview E/F columns were narrowed (the status text got shorter).
$ws1.Columns("E:F").ColumnWidth = 17.2159881591797

# --- zh-cn sheet -----------------------------------------------------------
$ws2.Range("C2:C3").Value = "Ready for handoff"
$ws2.Range("E2:E3").Value = "mt"
$ws2.Range("H2").Value = "2016-11-09 01:39:08"
$ws2.Range("H3").Value = "2016-11-09 01:39:08"
$ws2.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/121a1c3e8a5316928dfdff6199f4a99ecd944d89/e2e/5c36ebef-f17a-4cf5-a74d-5d65eb2bf5ec.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/7272d79ca4387b78d3635f0dc88fcc84eeb510e5/e2e/5c36ebef-f17a-4cf5-a74d-5d65eb2bf5ec.md."

$ws2.Columns("C:C").ColumnWidth = 17.2159881591797
$ws2.Columns("P:P").ColumnWidth = 40

# --- de-de sheet -----------------------------------------------------------
$ws3.Range("C2:C3").Value = "Ready for handoff"
$ws3.Range("E2:E3").Value = "mt"
$ws3.Range("H2").Value = "2016-11-09 01:39:24"
$ws3.Range("H3").Value = "2016-11-09 01:39:24"
$ws3.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/121a1c3e8a5316928dfdff6199f4a99ecd944d89/e2e/5c36ebef-f17a-4cf5-a74d-5d65eb2bf5ec.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/7272d79ca4387b78d3635f0dc88fcc84eeb510e5/e2e/5c36ebef-f17a-4cf5-a74d-5d65eb2bf5ec.md."

$ws3.Columns("C:C").ColumnWidth = 17.2159881591797
$ws3.Columns("P:P").ColumnWidth = 40
